# New Bill Payment Push
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (FEWA / partial) - account number changes to a new FEWA account
$ws.Range("F2").Value = "011074763011"

# Row 3 becomes a "Du" / "du Postpaid" bill, with its own consumer number and account number
$ws.Range("A3").Value = "Du"
$ws.Range("B3").Value = "du Postpaid"
$ws.Range("C3").Value = 599588538
$ws.Range("F3").Value = "011074763023"

# Row 4 (new) - ADDC bill push, re-using the first account number
$ws.Range("A4").Value = "ADDC"
$ws.Range("B4").Value = "ADDC"
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = 8235355929
$ws.Range("F4").Value = "011074763011"

# Update selection to reflect the new active cell after the push
[void]$ws.Range("F4").Select()
